# Work on main doors logic
# Insert a new "Cage Numbering System" slide before the "Old Slides" divider
# slide (position 6) and populate it with the diagram shapes.

$p = $ppt.ActivePresentation

# EMU -> points conversion helper (1 pt = 12700 EMU)
function Pt($emu) { return $emu / 12700 }

# --- Insert the new slide at position 6, using the Blank layout ------------
$blankLayout = $p.SlideMaster.CustomLayouts.Item(7)
$slide = $p.Slides.AddSlide(6, $blankLayout)

# --- Shape: Rectangle 3 (big translucent purple panel) ----------------------
$rect3 = $slide.Shapes.AddShape(1, (Pt 956929), (Pt 3476843), (Pt 6719777), (Pt 10866474))
$rect3.Name = "Rectangle 3"
$rect3.Fill.ForeColor.RGB = 0xC53F8C
$rect3.Fill.Transparency = 0.3
$rect3.Line.ForeColor.ObjectThemeColor = 1
$rect3.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle 5 (top-right translucent accent1 panel) --------------
$rect5 = $slide.Shapes.AddShape(1, (Pt 12408198), (Pt 2434851), (Pt 8580473), (Pt 4561367))
$rect5.Name = "Rectangle 5"
$rect5.Fill.ForeColor.ObjectThemeColor = 5
$rect5.Fill.Transparency = 0.3
$rect5.Line.ForeColor.ObjectThemeColor = 1
$rect5.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle 7 (bottom-right translucent accent1 panel) -----------
$rect7 = $slide.Shapes.AddShape(1, (Pt 12408198), (Pt 10823940), (Pt 8580473), (Pt 4561367))
$rect7.Name = "Rectangle 7"
$rect7.Fill.ForeColor.ObjectThemeColor = 5
$rect7.Fill.Transparency = 0.3
$rect7.Line.ForeColor.ObjectThemeColor = 1
$rect7.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: TextBox 8 ("2") --------------------------------------------------
$tb8 = $slide.Shapes.AddTextbox(1, (Pt 2583710), (Pt 8402248), (Pt 3466214), (Pt 1015663))
$tb8.Name = "TextBox 8"
$tb8.TextFrame.WordWrap = -1
$tb8.TextFrame.AutoSize = 1
$tb8.TextFrame.VerticalAnchor = 3
$tb8.TextFrame.TextRange.Text = "2"
$tb8.TextFrame.TextRange.Font.Size = 60
$tb8.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: TextBox 9 ("Cage Numbering System" title) -----------------------
$tb9 = $slide.Shapes.AddTextbox(1, (Pt 914400), (Pt 340242), (Pt 19436316), (Pt 1246495))
$tb9.Name = "TextBox 9"
$tb9.TextFrame.WordWrap = -1
$tb9.TextFrame.AutoSize = 1
$tb9.TextFrame.TextRange.Text = "Cage Numbering System"
$tb9.TextFrame.TextRange.Font.Size = 75

# --- Shape: Rectangle 4 (green box, top) ------------------------------------
$rect4 = $slide.Shapes.AddShape(1, (Pt 8173779), (Pt 3476843), (Pt 3721396), (Pt 2477386))
$rect4.Name = "Rectangle 4"
$rect4.Fill.ForeColor.RGB = 0x50B000
$rect4.Fill.Transparency = 0.4
$rect4.Line.ForeColor.ObjectThemeColor = 1
$rect4.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle: Rounded Corners 10 -----------------------------------
$rr10 = $slide.Shapes.AddShape(5, (Pt 8460863), (Pt 4040372), (Pt 1446028), (Pt 1339702))
$rr10.Name = "Rectangle: Rounded Corners 10"
$rr10.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle: Rounded Corners 11 -----------------------------------
$rr11 = $slide.Shapes.AddShape(5, (Pt 10162063), (Pt 4040372), (Pt 1446028), (Pt 1339702))
$rr11.Name = "Rectangle: Rounded Corners 11"
$rr11.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle 14 (green box, bottom) --------------------------------
$rect14 = $slide.Shapes.AddShape(1, (Pt 8173779), (Pt 11865930), (Pt 3721396), (Pt 2477386))
$rect14.Name = "Rectangle 14"
$rect14.Fill.ForeColor.RGB = 0x50B000
$rect14.Fill.Transparency = 0.4
$rect14.Line.ForeColor.ObjectThemeColor = 1
$rect14.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle: Rounded Corners 15 -----------------------------------
$rr15 = $slide.Shapes.AddShape(5, (Pt 8460863), (Pt 12429459), (Pt 1446028), (Pt 1339702))
$rr15.Name = "Rectangle: Rounded Corners 15"
$rr15.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape: Rectangle: Rounded Corners 16 -----------------------------------
$rr16 = $slide.Shapes.AddShape(5, (Pt 10162063), (Pt 12429459), (Pt 1446028), (Pt 1339702))
$rr16.Name = "Rectangle: Rounded Corners 16"
$rr16.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Connector 18 (flipped, above top green box) ----------------------------
$x = Pt 8173779
$y = Pt 3104707
$cx = Pt 3721396
$c18 = $slide.Shapes.AddConnector(1, ($x + $cx), $y, $x, $y)
$c18.Name = "Straight Arrow Connector 18"
$c18.HorizontalFlip = -1
$c18.Height = 0
$c18.Line.Weight = 5
$c18.Line.EndArrowheadStyle = 2

# --- Connector 20 (flipped, above bottom green box) -------------------------
$x = Pt 8173779
$y = Pt 11461897
$cx = Pt 3721396
$c20 = $slide.Shapes.AddConnector(1, ($x + $cx), $y, $x, $y)
$c20.Name = "Straight Arrow Connector 20"
$c20.HorizontalFlip = -1
$c20.Height = 0
$c20.Line.Weight = 5
$c20.Line.EndArrowheadStyle = 2

# --- Connector 22 (below top green box) -------------------------------------
$x = Pt 8173779
$y = Pt 6337005
$cx = Pt 3721396
$c22 = $slide.Shapes.AddConnector(1, $x, $y, ($x + $cx), $y)
$c22.Name = "Straight Arrow Connector 22"
$c22.Line.Weight = 5
$c22.Line.EndArrowheadStyle = 2

# --- Connector 23 (below bottom green box) ----------------------------------
$x = Pt 8173779
$y = Pt 14694196
$cx = Pt 3721396
$c23 = $slide.Shapes.AddConnector(1, $x, $y, ($x + $cx), $y)
$c23.Name = "Straight Arrow Connector 23"
$c23.Line.Weight = 5
$c23.Line.EndArrowheadStyle = 2
